$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(119).Insert()

$ws.Cells.Item(119, 1).Value = 8
$ws.Cells.Item(119, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44455
$ws.Cells.Item(119, 5).Value = 4
$ws.Cells.Item(119, 6).Value = 100112032
$ws.Cells.Item(119, 7).Value = "Zapallo italiano"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 600
$ws.Cells.Item(119, 11).Value = 10000
$ws.Cells.Item(119, 12).Value = 11000
$ws.Cells.Item(119, 13).Value = 10500
$ws.Cells.Item(119, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 210
$ws.Cells.Item(119, 17).Value = 50
$ws.Cells.Item(119, 18).Value = "Hortaliza"
